# "petites modifs code industrie" - rename the RESOURCE(S) related sheets
# and refresh the related title cell / sheet-view state to match.

$wb = $excel.ActiveWorkbook

# --- Rename sheets -----------------------------------------------------
# RESOURCE -> RESOURCES
$wsResources = $wb.Worksheets.Item("RESOURCE")
$wsResources.Name = "RESOURCES"

# TECHNOLOGIES_RESOURCE -> TECHNOLOGIES_RESOURCES
$wsTechResources = $wb.Worksheets.Item("TECHNOLOGIES_RESOURCE")
$wsTechResources.Name = "TECHNOLOGIES_RESOURCES"

$wsTechnologies = $wb.Worksheets.Item("TECHNOLOGIES")

# --- Update the sheet title cell to match the new sheet name -----------
$wsResources.Range("A1").Value = "RESOURCES"

# --- Refresh view / selection state -------------------------------------
# TECHNOLOGIES_RESOURCES: selection moved from C12 down to C18
$wsTechResources.Range("C18").Select()

# RESOURCES: drop the stale A2 selection back to the top-left cell
$wsResources.Range("A1").Select()

# TECHNOLOGIES becomes the active / selected tab (it was TECHNOLOGIES_RESOURCES before)
$wsTechnologies.Activate()
$wsTechnologies.Range("F1").Select()
